$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 74 (Serie index 84, 01-01-2021) updated values
$ws.Range("B74").Value = 22232
$ws.Range("C74").Value = 13109
$ws.Range("D74").Value = 12090
$ws.Range("E74").Value = 6218
$ws.Range("F74").Value = 2905
$ws.Range("G74").Value = 17967
$ws.Range("H74").Value = 5803
$ws.Range("J74").Value = 7075
$ws.Range("K74").Value = 2468

# Row 75 (Serie index 85, 01-04-2021) updated values
$ws.Range("B75").Value = 23234
$ws.Range("C75").Value = 14899
$ws.Range("D75").Value = 13636
$ws.Range("E75").Value = 6935
$ws.Range("G75").Value = 19060
$ws.Range("H75").Value = 6106
$ws.Range("J75").Value = 7611
$ws.Range("K75").Value = 2668
$ws.Range("L75").Value = 4287
